$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.583.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.649.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +10.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "579.13"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.644.89"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.681"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.63"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000296"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +18.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.14"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.229.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.641.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.482.45"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "405.20"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +28.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.16"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.97"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.70"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.36%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +9.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +25.38%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "696.28"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +16.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.31"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.97"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.96"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.431"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +16.74%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0796"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.143"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +11.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.93"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +22.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +14.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.209.91"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +11.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.02"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +37.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.97"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.73%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.45%  "
